# Applies the "added report last 11-11-24" update to the RSO Monthly
# Incentive Sheet (Oct-2024) — updates GA achievement for Asim Gain,
# updates airtime achievement figures for three RSOs, and records the
# newly-added "Eid Bonus" (column L) of 1000 for every RSO row.
# All dependent formulas (G, K, M, P columns and the row-8 / row-11
# subtotal rows, as well as the grand total on row 16) recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 4 : Asim Gain ---------------------------------------------------
$ws.Range("F4").Value = 41
$ws.Range("L4").Value = 1000

# --- Row 5 : Ismil Hossain ------------------------------------------------
$ws.Range("J5").Value = 1124359
$ws.Range("L5").Value = 1000

# --- Row 6 : Sukanto Sarkar ------------------------------------------------
$ws.Range("J6").Value = 1241020
$ws.Range("L6").Value = 1000

# --- Row 7 : Md Ashikur Rahman --------------------------------------------
$ws.Range("L7").Value = 1000

# --- Update the scroll position / selection left by the author -----------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("G24").Select()
